$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K19").Value = 2139.4285
$ws.Range("M19").Value = -1964.4285
$ws.Range("I19").Value = 2139.4285
$ws.Range("H19").Value = 2647.9167
$ws.Range("M28").Value = -105.9091
$ws.Range("H28").Value = 926.6875
$ws.Range("I28").Value = 590.9091
$ws.Range("K28").Value = 590.9091
$ws.Range("L88").Value = 45000
$ws.Range("H88").Value = 45000
$ws.Range("N88").Value = -45812
$ws.Range("J88").Value = 45000
$ws.Range("J91").Value = 45000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47808
$ws.Range("H91").Value = 45000
$ws.Range("N99").Value = -15118.4
$ws.Range("L99").Value = 12122.4
$ws.Range("H99").Value = 3622.7273
$ws.Range("J99").Value = 4040.8
$ws.Range("N112").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("I116").Value = 87995.164
$ws.Range("M116").Value = -84553.164
$ws.Range("K116").Value = 87995.164
$ws.Range("H116").Value = 53891.332
$ws.Range("H132").Value = 47303.273
$ws.Range("I132").Value = 49488.855
$ws.Range("M132").Value = -145936.565
$ws.Range("K132").Value = 148466.565
$ws.Range("K138").Value = 11222.7855
$ws.Range("I138").Value = 3740.9285
$ws.Range("J138").Value = 8299.451999999999
$ws.Range("L138").Value = 24898.356
$ws.Range("N138").Value = -35178.356
$ws.Range("H138").Value = 7159.8213
$ws.Range("M138").Value = -6082.7855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 8153.7334
$ws.Range("H2").Value = 6926.6523
$ws.Range("K2").Value = 8153.7334
$ws.Range("M2").Value = -8040.7334
$ws.Range("H5").Value = 343.25
$ws.Range("I5").Value = 87.5
$ws.Range("M5").Value = 24.5
$ws.Range("K5").Value = 87.5
$ws.Range("I32").Value = 3294.1333
$ws.Range("K32").Value = 3294.1333
$ws.Range("H32").Value = 3248.3225
$ws.Range("M32").Value = -3007.1333
$ws.Range("N96").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("I116").Value = 8153.7334
$ws.Range("M116").Value = -5859.7334
$ws.Range("K116").Value = 8153.7334
$ws.Range("H116").Value = 6926.6523

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6926.6523
$ws.Range("K3").Value = 8153.7334
$ws.Range("I3").Value = 8153.7334
$ws.Range("M3").Value = -8039.7334
$ws.Range("H4").Value = 343.25
$ws.Range("K4").Value = 87.5
$ws.Range("M4").Value = 27.5
$ws.Range("I4").Value = 87.5
$ws.Range("N20").Value = -5422
$ws.Range("L20").Value = 4928
$ws.Range("M20").Value = -1936.5
$ws.Range("H20").Value = 3431
$ws.Range("I20").Value = 2183.5
$ws.Range("J20").Value = 4928
$ws.Range("K20").Value = 2183.5
$ws.Range("M24").ClearContents()
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 20018
$ws.Range("I24").Value = 0
$ws.Range("H24").Value = 20018
$ws.Range("N24").Value = -20488
$ws.Range("J24").Value = 20018
$ws.Range("I99").Value = 1241.8334
$ws.Range("H99").Value = 1241.8334
$ws.Range("K99").Value = 1241.8334
$ws.Range("M99").Value = 256.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 2859.9
$ws.Range("H31").Value = 2317.6
$ws.Range("I31").Value = 2859.9
$ws.Range("M31").Value = -2564.9
$ws.Range("K34").Value = 2859.9
$ws.Range("M34").Value = -2657.9
$ws.Range("H34").Value = 2317.6
$ws.Range("I34").Value = 2859.9
$ws.Range("H58").Value = 2189.7856
$ws.Range("M58").Value = -1852.6667
$ws.Range("K58").Value = 2055.6667
$ws.Range("I58").Value = 2055.6667
$ws.Range("I134").Value = 535.1
$ws.Range("K134").Value = 1605.3
$ws.Range("M134").Value = 929.6999999999998
$ws.Range("H134").Value = 535.1
$ws.Range("K136").Value = 6167.000100000001
$ws.Range("I136").Value = 2055.6667
$ws.Range("M136").Value = -3617.000100000001
$ws.Range("H136").Value = 2189.7856

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1780.909
$ws.Range("N107").Value = -10172.1666
$ws.Range("L107").Value = 6332.1666
$ws.Range("J107").Value = 2110.7222
$ws.Range("M110").ClearContents()
$ws.Range("K110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("H110").Value = 0
$ws.Range("K131").Value = 4166.25
$ws.Range("M131").Value = 873.75
$ws.Range("H131").Value = 112435.78
$ws.Range("I131").Value = 1388.75
$ws.Range("L131").Value = 603820.2
$ws.Range("N131").Value = -613900.2
$ws.Range("J131").Value = 201273.4
$ws.Range("H137").Value = 752229.8
$ws.Range("N137").Value = -3231459.6
$ws.Range("M137").Value = -925.0002000000004
$ws.Range("K137").Value = 6025.0002
$ws.Range("I137").Value = 2008.3334
$ws.Range("J137").Value = 1073753.2
$ws.Range("L137").Value = 3221259.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K11").Value = 3663156.8
$ws.Range("M11").Value = -3663017.8
$ws.Range("L11").Value = 8642857
$ws.Range("N11").Value = -8643135
$ws.Range("H11").Value = 5323057
$ws.Range("I11").Value = 3663156.8
$ws.Range("J11").Value = 8642857
$ws.Range("I43").Value = 837000.3
$ws.Range("K43").Value = 837000.3
$ws.Range("M43").Value = -836849.3
$ws.Range("H43").Value = 837000.3
$ws.Range("N46").Value = -30358
$ws.Range("I46").Value = 22000
$ws.Range("L46").Value = 30046
$ws.Range("J46").Value = 30046
$ws.Range("K46").Value = 22000
$ws.Range("M46").Value = -21844
$ws.Range("H46").Value = 24682
$ws.Range("M80").ClearContents()
$ws.Range("K80").Value = 0
$ws.Range("H80").Value = 1494.6666
$ws.Range("I80").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("K83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("H83").Value = 1494.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K22").Value = 792.55554
$ws.Range("M22").Value = -497.55554
$ws.Range("I22").Value = 792.55554
$ws.Range("H22").Value = 4540.2
$ws.Range("K27").Value = 792.55554
$ws.Range("H27").Value = 4540.2
$ws.Range("M27").Value = -685.55554
$ws.Range("I27").Value = 792.55554
$ws.Range("M30").ClearContents()
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("J31").Value = 8666.333000000001
$ws.Range("K31").Value = 1895.5
$ws.Range("I31").Value = 1895.5
$ws.Range("L31").Value = 8666.333000000001
$ws.Range("H31").Value = 4797.2856
$ws.Range("M31").Value = -1647.5
$ws.Range("N31").Value = -9162.333000000001
$ws.Range("I40").Value = 3382
$ws.Range("M40").Value = -3246
$ws.Range("K40").Value = 3382
$ws.Range("H40").Value = 5555.8887
$ws.Range("I122").Value = 3600.75
$ws.Range("H122").Value = 4132.3335
$ws.Range("M122").Value = -8352.25
$ws.Range("K122").Value = 10802.25
$ws.Range("H132").Value = 3701.8
$ws.Range("L132").Value = 11250
$ws.Range("N132").Value = -16310
$ws.Range("J132").Value = 3750
$ws.Range("K136").Value = 13195.5
$ws.Range("I136").Value = 4398.5
$ws.Range("M136").Value = -10645.5
$ws.Range("H136").Value = 125004720

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N126").Value = -15078.5
$ws.Range("H126").Value = 3215.8333
$ws.Range("L126").Value = 10138.5
$ws.Range("J126").Value = 3379.5
$ws.Range("I132").Value = 4931.067
$ws.Range("H132").Value = 5223.1577
$ws.Range("M132").Value = -12263.201
$ws.Range("J132").Value = 6318.5
$ws.Range("N132").Value = -24015.5
$ws.Range("L132").Value = 18955.5
$ws.Range("K132").Value = 14793.201
$ws.Range("N136").ClearContents()
$ws.Range("J136").Value = 0
$ws.Range("M136").Value = 1050
$ws.Range("H136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("I136").Value = 500
$ws.Range("L136").Value = 0
